$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.301.22'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.687.53'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '679.90'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.11'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.60%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -1.29%  '
$ws.Range("E9").Value = '  -1.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.13'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.86%  '
$ws.Range("E11").Value = '  -1.68%  '
$ws.Range("E12").Value = '  -3.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.309.92'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.47'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.692.52'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.330.85'
$ws.Range("D16").ClearFormats()
$ws.Range("E17").Value = '  +1.54%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '16.02'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.87%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.42'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '467.88'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.01%  '
$ws.Range("E21").Value = '  -0.50%  '
$ws.Range("E22").Value = '  -2.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '79.81'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.835.62'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.39%  '
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("E26").Value = '  -6.80%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.95'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.15'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -4.22%  '
$ws.Range("E29").Value = '  -2.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.75'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -4.31%  '
$ws.Range("E31").Value = '  -3.47%  '
$ws.Range("E32").Value = '  -4.49%  '
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.93'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.676.48'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.158'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -4.87%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.28'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.79%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.22'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.14%  '
$ws.Range("E40").Value = '  -2.76%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("E42").Value = '  -2.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '171.44'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +4.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.944'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '47.63'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.98%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.39'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -6.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.11'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.96%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.70'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -4.67%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000276'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -4.01%  '
$ws.Range("E50").Value = '  -5.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.79'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -3.96%  '
